$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 101, shifting rows 101:118 down to 102:119.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record's data.
$ws.Cells.Item(101, 1).Value = 4
$ws.Cells.Item(101, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(101, 3).Value = "Los Lagos"
$ws.Cells.Item(101, 4).Value = 44504
$ws.Cells.Item(101, 5).Value = 10
$ws.Cells.Item(101, 6).Value = 100112009
$ws.Cells.Item(101, 7).Value = "Acelga"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 100
$ws.Cells.Item(101, 11).Value = 3000
$ws.Cells.Item(101, 12).Value = 3000
$ws.Cells.Item(101, 13).Value = 3000
$ws.Cells.Item(101, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(101, 15).Value = "Región del Maule"
$ws.Cells.Item(101, 16).Value = 750
$ws.Cells.Item(101, 17).Value = 4
$ws.Cells.Item(101, 18).Value = "Hortaliza"

# Ensure the date cell carries the same custom date/time style as the other
# rows in this column (style index 2 in the original workbook).
$ws.Cells.Item(101, 4).NumberFormat = $ws.Cells.Item(102, 4).NumberFormat
